# Reposition two shapes on the single slide.
# PowerPoint's Shape.Left/.Top are Single-precision points (1 pt = 12700 EMU);
# the literals below are the closest float32 point values whose EMU
# round-trip (via truncation) lands exactly on the target EMU coordinates.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Shape 224" (text "5K"): keep Left (7924997 EMU), move Top
# from 3999765 EMU to 3960009 EMU.
$shape224 = $s.Shapes.Item("Shape 224")
$shape224.Left = 624.0155639648438
$shape224.Top = 311.8117370605469

# "Shape 334" (text "Crystal oscillator"): move from
# (6055602, 2518508) EMU to (6444181, 2027978) EMU.
$shape334 = $s.Shapes.Item("Shape 334")
$shape334.Left = 507.41583251953125
$shape334.Top = 159.68331909179688
